# Initial update with CPL's work-to-date
#
# Core content change: on the "HPEbP" sheet, the electrolysis efficiency
# formula in B3 drops the old "+46" term from its denominator
# (118/(162+2+46) -> 118/(162+2)); the shared-formula chain across the
# row (C3:AI3) recalculates from that new value. The workbook is also
# switched to iterative calculation, and the saved UI state (active
# sheet/selection on each tab) is updated to match.

$wb = $excel.ActiveWorkbook

# --- Calculation options: enable iterative calculation -----------------
$excel.Iteration = $true
$excel.MaxChange = 0.00001

$wsAbout = $wb.Worksheets.Item("About")
$wsIEA   = $wb.Worksheets.Item("IEA Data")
$wsHPEbP = $wb.Worksheets.Item("HPEbP")

# --- Core formula edit ---------------------------------------------------
$wsHPEbP.Range("B3").Formula = "=118/(162+2)"

# --- Selection / active-tab bookkeeping (matches the saved UI state) -----
$wsAbout.Range("B14").Select()

$wsIEA.Range("D7:F7").Select()

$wsHPEbP.Activate()
$wsHPEbP.Range("C3").Select()
